$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing data block (rows 203-238) down by one row (to rows 204-239),
# preserving values, types and styles.
$src = $ws.Range("A203:R238")
$dst = $ws.Range("A204:R239")
$src.Copy($dst)

# Write the new record into row 203 (a new weekly price observation).
$ws.Range("D203").Value = 44504
$ws.Range("J203").Value = 250
$ws.Range("K203").Value = 11000
$ws.Range("L203").Value = 11000
$ws.Range("M203").Value = 11000
$ws.Range("P203").Value = 550

Write-Host "done"
